$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "FE-0317000"
$ws.Range("C2").Value = "07/08/2024 11:42:00"
$ws.Range("D2").Value = "07/08/2024 13:05:00"
$ws.Range("G2").Value = 4980
$ws.Range("H2").Value = 3000

# Update row 3 values
$ws.Range("A3").Value = "FE-0317001"
$ws.Range("C3").Value = "07/08/2024 12:47:00"
$ws.Range("D3").Value = "07/08/2024 13:05:00"
$ws.Range("G3").Value = 1080

# Remove rows 4 and 5 (no longer part of the exported data)
$ws.Range("A4:H5").EntireRow.Delete()
